# Add the new "simple_charcoal" unit process as row 47 of the
# "Unit Processes" sheet, matching the columns used by the existing rows
# (e.g. row 46 "simple_casting" / row 45 "simple_syngas").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Unit Processes")

$ws.Range("A47").Value = "simple_charcoal"
$ws.Range("B47").Value = "energy"
$ws.Range("C47").Value = "Charcoal"
$ws.Range("D47").Value = "charcoal"
$ws.Range("E47").Value = "outflow"
$ws.Range("F47").Value = "data/steel/steel_simplified_var.xlsx"
$ws.Range("G47").Value = "Charcoal"
$ws.Range("H47").Value = "data/steel/steel_simplified_calcs.xlsx"
$ws.Range("I47").Value = "Charcoal"

# Copy the formatting (text number format etc.) from the row above so the
# new row's cells pick up the same style as the rest of the table.
$ws.Range("A46:I46").Copy()
$ws.Range("A47:I47").PasteSpecial(-4122)

# Match the updated selection recorded in the saved file.
$ws.Range("I47").Select()
